$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.458.03"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.088.18"
$ws.Range("E3").Value = "  +4.26%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.55"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.27"
$ws.Range("E6").Value = "  +3.21%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  +6.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.086.63"
$ws.Range("E10").Value = "  +4.21%  "
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("E13").Value = "  +6.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.609.70"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.01"
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.439.22"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000194"
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.064.98"
$ws.Range("E18").Value = "  +3.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.62"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.09"
$ws.Range("E20").Value = "  +4.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.49"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.47"
$ws.Range("E22").Value = "  +9.24%  "
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.230.80"
$ws.Range("E24").Value = "  +3.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.42"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.34"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.88"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "501.31"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  +12.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.72"
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.76"
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.06"
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "192.91"
$ws.Range("E40").Value = "  +6.67%  "
$ws.Range("E41").Value = "  -4.59%  "
$ws.Range("E42").Value = "  -9.26%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.793"
$ws.Range("E44").Value = "  +19.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.12"
$ws.Range("E45").Value = "  +3.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.25"
$ws.Range("E46").Value = "  +4.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.25"
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("E49").Value = "  +4.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.598"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.88"
$ws.Range("E51").Value = "  -0.48%  "
